# Duplicate the most recent dated sheet ("02-17-21") to create a new
# "08-01-22" sheet immediately after it, keeping all of its data,
# formatting (pink header fill) and page setup intact, then make sure
# the header row still reads A / B / C and leave the new sheet as the
# active tab (mirrors a user right-clicking the tab -> Move or Copy...
# -> Create a copy, then renaming it).

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("02-17-21")
$template.Copy([System.Reflection.Missing]::Value, $template)

$newSheet = $wb.Worksheets.Item($template.Index + 1)
$newSheet.Name = "08-01-22"

# Keep the header labels as A / B / C on the new sheet.
$newSheet.Range("A1").Value = "A"
$newSheet.Range("B1").Value = "B"
$newSheet.Range("C1").Value = "C"

# The newly copied/renamed sheet is the one the user is left looking at.
$newSheet.Activate()
